$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "consumer/mad/runme_large.sh"
$ws.Range("B8").Value = 0.03
$ws.Range("C8").Value = 0.02
$ws.Range("D8").Value = 0.01

$ws.Range("A10").Select() | Out-Null
